# EN COURS - task references #133: Creation ip_temporaire
# Adds the 5 new "ip_temporaire" related queries as rows 25-29 to the
# "liste des requetes" sheet (columns: Table concernée, Type, Nom requete,
# Fonctionnalités, Requete).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$table      = @("ip_temporaire", "ip_temporaire", "ip_temporaire", "ip_temporaire", "ip_temporaire")
$type       = @("Insert", "Update", "Select", "Select", "Delete")
$nomRequete = @(
    "REALMS_INS_IPTEMPORAIRE_STOCKAGEIPTEMPORAIRE",
    "REALMS_UPD_IPTEMPORAIRE_MAJIPTEMPORAIRE",
    "REALMS_SEL_IPTEMPORAIRE_LECTURENERREURS",
    "REALMS_SEL_IPTEMPORAIRE_RECHERCHEIP",
    "REALMS_DEL_IPTEMPORAIRE_SUPPRLIGNEIP"
)
$fonction = @(
    "Création d'une erreur d'authentification sur une nouvelle ip",
    "Création d'une erreur d'authentification sur une ip existante",
    "Lecteure du nombre d'erreurs d'authentification sur une ip",
    "Savoir si une ip a déjà eu une erreur d'authentification",
    "Suppression des lignes sans erreurs d'authentification"
)
$requete = @(
    "INSERT INTO ip_temporaire VALUES (`$1, '1')",
    "UPDATE ip_temporaire SET ip_temp_nessais = `$1 WHERE ip_temp_ip = `$2`"",
    "SELECT ip_temp_nessais FROM ip_temporaire where ip_temp_ip = `$1",
    "SELECT ip_temp_ip FROM ip_temporaire where ip_temp_ip = `$1",
    "DELETE FROM ip_temporaire where ip_temp_nessais = '0'"
)

$startRow = 25

# Fill column by column (C, A, B, D, E) to reproduce the shared-string
# insertion order recorded in the source workbook.
for ($i = 0; $i -lt $nomRequete.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $nomRequete[$i]
}
for ($i = 0; $i -lt $table.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $table[$i]
}
for ($i = 0; $i -lt $type.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $type[$i]
}
for ($i = 0; $i -lt $fonction.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $fonction[$i]
}
for ($i = 0; $i -lt $requete.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $requete[$i]
}

# Update selection/view to mirror the recorded end-state: the window had
# scrolled so column D was the left-most visible column, row 7 the topmost
# visible row, with E30 as the active cell/selection.
$ws.Range("E30").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 4
